$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number (45203 = 2023-10-04)
# that must be bumped to 45204 (2023-10-05) for every data row (2..494).
$ws.Range("C2:C494").Value2 = 45204
